# MS549 Assignment 2 - Video and ZIP
# Update the "Delete/Remove -> .NET LinkedList" performance numbers (column F)
# on Sheet1, and move the active-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F3").Value = 27.21
$ws.Range("F4").Value = 1171.3699999999999
$ws.Range("F5").Value = 116726.36

$ws.Range("J6").Select()
